# Updates crypto price/volume values to match latest scrape

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.901.87"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.888.31"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7723"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.01"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.60"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07183"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08576"
$ws.Range("E11").Value = "  +5.95%  "
$ws.Range("D12").Value = "1.946.33"
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7638"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.376"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.75"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.194"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "29.917.89"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.36"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007813"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "2.176.21"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9979"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.973"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1647"
$ws.Range("E25").Value = "  +4.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.375"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.74"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.444"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.536"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.521"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.105"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05441"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.242"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7455"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.692"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01967"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4462"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").Value = "1.108.63"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.34"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.086"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.85"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.627"
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("D50").Value = "2.073.10"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.986"
$ws.Range("E51").Value = "  -0.31%  "
